$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
for ($i=0; $i -le 14; $i++) {
  $addr = "A" + (300+$i)
  $c = $ws.Range($addr)
  $b = $c.Borders.Item(7)
  $b.LineStyle = 1
  $b.ThemeColor = $i
  Write-Host $i "->" $b.ThemeColor
}
